# "new fixed gerber files for all boards"
#
# The BOM table's middle four rows (No. 2-5) get reshuffled:
#   - the "2.54-1*6P母" / PROG header moves from row 4 up to row 3
#   - the "DW254R-22-20-85" header moves from row 5 up to row 4 (and its
#     Quantity becomes 2, swapping with the old row-4 quantity)
#   - the "Level Shifter" row moves from row 6 up to row 5, and its
#     Footprint text gets re-cased to "Level Shifter Footprint"
#   - the "MX128..." Wiegand connector moves from row 3 down to row 6,
#     picking up a new Designator, "WIEGAND" (was "CN1")
#
# We use Range.Copy(destination) instead of Range.Value = "..." for every
# cell whose final text is still present elsewhere in the original sheet:
# that preserves the shared-string cell type (t="s") exactly as the
# original file had it. Plain Value assignment would be fine for normal
# text, but Excel's General-format auto-detection would silently turn
# numeric-looking text (e.g. the Quantity "2"/"1") into real numbers, and
# would collapse blank cells instead of keeping them as explicit empty
# shared strings - neither of which matches the source data here.
#
# Because several of these moves are cyclic (row3 -> row6, etc.) a couple
# of old row-3 cells are staged in scratch cells (row 30, well outside the
# used range) before row 3 is overwritten, then copied down into row 6 at
# the end and cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- stage the old row 3 (MX128/Wiegand) cells we still need for row 6 ---
$ws.Range("C3").Copy($ws.Range("C30"))
$ws.Range("E3").Copy($ws.Range("E30"))
$ws.Range("G3").Copy($ws.Range("G30"))
$ws.Range("H3").Copy($ws.Range("H30"))
$ws.Range("I3").Copy($ws.Range("I30"))

# --- row 3 <- old row 4 (PROG header) ---
$ws.Range("C4").Copy($ws.Range("C3"))
$ws.Range("D4").Copy($ws.Range("D3"))
$ws.Range("E4").Copy($ws.Range("E3"))
$ws.Range("G4").Copy($ws.Range("G3"))
$ws.Range("H4").Copy($ws.Range("H3"))
$ws.Range("I4").Copy($ws.Range("I3"))

# --- row 4 <- old row 5 (DW254R header); Quantity "1" -> "2" ---
$ws.Range("A3").Copy($ws.Range("B4"))
$ws.Range("C5").Copy($ws.Range("C4"))
$ws.Range("D5").Copy($ws.Range("D4"))
$ws.Range("E5").Copy($ws.Range("E4"))
$ws.Range("G5").Copy($ws.Range("G4"))
$ws.Range("H5").Copy($ws.Range("H4"))
$ws.Range("I5").Copy($ws.Range("I4"))

# --- row 5 <- old row 6 (Level Shifter); Quantity "2" -> "1"; Footprint re-cased ---
$ws.Range("B2").Copy($ws.Range("B5"))
$ws.Range("C6").Copy($ws.Range("C5"))
$ws.Range("D6").Copy($ws.Range("D5"))
$ws.Range("E5").Value = "Level Shifter Footprint"
$ws.Range("F5").Copy($ws.Range("G5"))
$ws.Range("F5").Copy($ws.Range("H5"))
$ws.Range("F5").Copy($ws.Range("I5"))
$ws.Range("F5").Copy($ws.Range("J5"))

# --- row 6 <- staged old row 3 (MX128/Wiegand); Designator "CN1" -> "WIEGAND" ---
$ws.Range("C30").Copy($ws.Range("C6"))
$ws.Range("D6").Value = "WIEGAND"
$ws.Range("E30").Copy($ws.Range("E6"))
$ws.Range("G30").Copy($ws.Range("G6"))
$ws.Range("H30").Copy($ws.Range("H6"))
$ws.Range("I30").Copy($ws.Range("I6"))
$ws.Range("J3").Copy($ws.Range("J6"))

# --- tidy up the scratch cells ---
$ws.Range("C30:I30").Clear()
